$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows' aspect_ratio (column I) values ---
$ws.Range("I2").Value = "4:5"
$ws.Range("I3").Value = "4:5"
$ws.Range("I4").Value = "16:9"

# --- Update likes_count / comments_count on row 3 ---
$ws.Range("T3").Value = 554
$ws.Range("U3").Value = 3

# --- Append new rows 5, 6, 7 of content data ---

# Row 5
$ws.Range("A5").Value = "the_year_book_"
$ws.Range("B5").Value = "Post"
$ws.Range("C5").Value = "DJecytozKis"
$ws.Range("D5").Value = "2025-05-10T14:06:03.000Z"
$ws.Range("E5").Value = "Saturday"
$ws.Range("F5").Value = 19
$ws.Range("G5").Value = "07:36 PM"
$ws.Range("I5").Value = "4:5"
$ws.Range("J5").Value = "Do bhai aur Bhagwan ki Kripa se dono tabahi😎🧿💎"
$ws.Range("K5").Value = 46
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = "[]"
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = "[]"
$ws.Range("P5").Value = 0
$ws.Range("T5").Value = 700
$ws.Range("U5").Value = 33
$ws.Range("V5").Value = "Adiyogi Coimbatore"

# Row 6
$ws.Range("A6").Value = "the_year_book_"
$ws.Range("B6").Value = "Post"
$ws.Range("C6").Value = "DJHpfpGSAcw"
$ws.Range("D6").Value = "2025-05-01T17:34:31.000Z"
$ws.Range("E6").Value = "Thursday"
$ws.Range("F6").Value = 23
$ws.Range("G6").Value = "11:04 PM"
$ws.Range("I6").Value = "4:5"
$ws.Range("J6").Value = "Om Namah Shivaya 🙏  #new #shiv #adiyogi"
$ws.Range("K6").Value = 40
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = "['#new', '#shiv', '#adiyogi']"
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = "[]"
$ws.Range("P6").Value = 0
$ws.Range("T6").Value = 199
$ws.Range("U6").Value = 7
$ws.Range("V6").Value = "Adiyogi Shiva statue"

# Row 7
$ws.Range("A7").Value = "the_year_book_"
$ws.Range("B7").Value = "Post"
$ws.Range("C7").Value = "DJB1QwfPUKC"
$ws.Range("D7").Value = "2025-04-29T11:21:53.000Z"
$ws.Range("E7").Value = "Tuesday"
$ws.Range("F7").Value = 16
$ws.Range("G7").Value = "04:51 PM"
$ws.Range("I7").Value = "4:5"
$ws.Range("J7").Value = "Eyes on me!👀 #new #post #blackoutfit"
$ws.Range("K7").Value = 37
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = "['#new', '#post', '#blackoutfit']"
$ws.Range("N7").Value = 3
$ws.Range("O7").Value = "[]"
$ws.Range("P7").Value = 0
$ws.Range("T7").Value = 114
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = "Barbil Keonjhar"

# --- Resize the table / autofilter to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:V7"))
